$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing I2:I5 values from 4 to 5
$ws.Range("I2").Value = 5
$ws.Range("I3").Value = 5
$ws.Range("I4").Value = 5
$ws.Range("I5").Value = 5

# Add new row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = -5
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 12
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = "train_dim2_1"

# Update selection to A6
$ws.Range("A6").Select()

$wb.Save()
